$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 3205566
$ws.Range("I8").Value = 3205566
$ws.Range("K8").Value = 9616698
$ws.Range("M8").Value = -9616559
$ws.Range("H17").Value = 10978.179
$ws.Range("J17").Value = 11418.808
$ws.Range("L17").Value = 34256.424
$ws.Range("N17").Value = -34592.424
$ws.Range("H51").Value = 4166.5
$ws.Range("I51").Value = 4799.9
$ws.Range("K51").Value = 4799.9
$ws.Range("M51").Value = -4315.9
$ws.Range("H96").Value = 1936.1538
$ws.Range("I96").Value = 1626.3636
$ws.Range("J96").Value = 3640
$ws.Range("K96").Value = 4879.0908
$ws.Range("L96").Value = 10920
$ws.Range("M96").Value = -3506.0908
$ws.Range("N96").Value = -13666
$ws.Range("H98").Value = 1270.9131
$ws.Range("I98").Value = 1204.3158
$ws.Range("J98").Value = 1587.25
$ws.Range("K98").Value = 1204.3158
$ws.Range("L98").Value = 1587.25
$ws.Range("M98").Value = 293.6841999999999
$ws.Range("N98").Value = -4583.25
$ws.Range("H122").Value = 1270.9131
$ws.Range("I122").Value = 1204.3158
$ws.Range("J122").Value = 1587.25
$ws.Range("K122").Value = 3612.9474
$ws.Range("L122").Value = 4761.75
$ws.Range("M122").Value = -1162.9474
$ws.Range("N122").Value = -9661.75
$ws.Range("H132").Value = 3323.037
$ws.Range("I132").Value = 3348
$ws.Range("K132").Value = 10044
$ws.Range("M132").Value = -7514
$ws.Range("H138").Value = 3045.8367
$ws.Range("I138").Value = 1138.3334
$ws.Range("K138").Value = 3415.0002
$ws.Range("M138").Value = 1724.9998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11589459
$ws.Range("I32").Value = 11971159
$ws.Range("K32").Value = 11971159
$ws.Range("M32").Value = -11970872
$ws.Range("H61").Value = 3664.1892
$ws.Range("I61").Value = 3864.963
$ws.Range("J61").Value = 3122.1
$ws.Range("K61").Value = 3864.963
$ws.Range("L61").Value = 3122.1
$ws.Range("M61").Value = -3652.963
$ws.Range("N61").Value = -3546.1
$ws.Range("H63").Value = 4914.9
$ws.Range("J63").Value = 5529.8
$ws.Range("L63").Value = 5529.8
$ws.Range("N63").Value = -6901.8
$ws.Range("H66").Value = 4914.9
$ws.Range("J66").Value = 5529.8
$ws.Range("L66").Value = 27649
$ws.Range("N66").Value = -34513
$ws.Range("H132").Value = 2708.8408
$ws.Range("I132").Value = 2240.647
$ws.Range("J132").Value = 4300.7
$ws.Range("K132").Value = 6721.941
$ws.Range("L132").Value = 12902.1
$ws.Range("M132").Value = -4191.941
$ws.Range("N132").Value = -17962.1
$ws.Range("H136").Value = 3664.1892
$ws.Range("I136").Value = 3864.963
$ws.Range("J136").Value = 3122.1
$ws.Range("K136").Value = 11594.889
$ws.Range("L136").Value = 9366.299999999999
$ws.Range("M136").Value = -9044.889000000001
$ws.Range("N136").Value = -14466.3

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2876.76
$ws.Range("I99").Value = 1612.4615
$ws.Range("K99").Value = 1612.4615
$ws.Range("M99").Value = -114.4614999999999
$ws.Range("H105").Value = 2400.8
$ws.Range("I105").Value = 2147.5264
$ws.Range("J105").Value = 3202.8333
$ws.Range("K105").Value = 2147.5264
$ws.Range("L105").Value = 3202.8333
$ws.Range("M105").Value = -400.5264000000002
$ws.Range("N105").Value = -6696.8333
$ws.Range("H134").Value = 1702787.9
$ws.Range("I134").Value = 2646641.5
$ws.Range("J134").Value = 3851.1333
$ws.Range("K134").Value = 7939924.5
$ws.Range("L134").Value = 11553.3999
$ws.Range("M134").Value = -7937389.5
$ws.Range("N134").Value = -16623.3999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 2399.6
$ws.Range("I2").Value = 2399.6
$ws.Range("K2").Value = 2399.6
$ws.Range("M2").Value = -2286.6
$ws.Range("H31").Value = 3162.6858
$ws.Range("I31").Value = 2195.35
$ws.Range("J31").Value = 4452.467
$ws.Range("K31").Value = 2195.35
$ws.Range("L31").Value = 4452.467
$ws.Range("M31").Value = -1900.35
$ws.Range("N31").Value = -5042.467
$ws.Range("H34").Value = 3162.6858
$ws.Range("I34").Value = 2195.35
$ws.Range("J34").Value = 4452.467
$ws.Range("K34").Value = 2195.35
$ws.Range("L34").Value = 4452.467
$ws.Range("M34").Value = -1993.35
$ws.Range("N34").Value = -4856.467
$ws.Range("H58").Value = 2892.1936
$ws.Range("I58").Value = 2736.08
$ws.Range("K58").Value = 2736.08
$ws.Range("M58").Value = -2533.08
$ws.Range("H122").Value = 9095265
$ws.Range("I122").Value = 11115824
$ws.Range("K122").Value = 33347472
$ws.Range("M122").Value = -33345022
$ws.Range("H134").Value = 2480.9666
$ws.Range("I134").Value = 2190.2
$ws.Range("J134").Value = 3062.5
$ws.Range("K134").Value = 6570.599999999999
$ws.Range("L134").Value = 9187.5
$ws.Range("M134").Value = -4035.599999999999
$ws.Range("N134").Value = -14257.5
$ws.Range("H136").Value = 2892.1936
$ws.Range("I136").Value = 2736.08
$ws.Range("K136").Value = 8208.24
$ws.Range("M136").Value = -5658.24
$ws.Range("H139").Value = 90000
$ws.Range("J139").Value = 90000
$ws.Range("L139").Value = 90000
$ws.Range("N139").Value = -100280

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 2299
$ws.Range("J129").Value = 3162.8
$ws.Range("L129").Value = 9488.400000000001
$ws.Range("N129").Value = -19488.4
$ws.Range("H131").Value = 1614.5454
$ws.Range("I131").Value = 923.1667
$ws.Range("J131").Value = 1873.8125
$ws.Range("K131").Value = 2769.5001
$ws.Range("L131").Value = 5621.4375
$ws.Range("M131").Value = 2270.4999
$ws.Range("N131").Value = -15701.4375
$ws.Range("H132").Value = 1257.6
$ws.Range("I132").Value = 1466.3334
$ws.Range("J132").Value = 944.5
$ws.Range("K132").Value = 13197.0006
$ws.Range("L132").Value = 8500.5
$ws.Range("M132").Value = -10667.0006
$ws.Range("N132").Value = -13560.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 103999
$ws.Range("J15").Value = 103999
$ws.Range("L15").Value = 103999
$ws.Range("N15").Value = -104575
$ws.Range("H81").Value = 103999
$ws.Range("J81").Value = 103999
$ws.Range("L81").Value = 103999
$ws.Range("N81").Value = -105995
$ws.Range("H84").Value = 103999
$ws.Range("J84").Value = 103999
$ws.Range("L84").Value = 311997
$ws.Range("N84").Value = -321981
$ws.Range("H97").Value = 973.35
$ws.Range("I97").Value = 869.7778
$ws.Range("K97").Value = 869.7778
$ws.Range("M97").Value = -373.7778
$ws.Range("H122").Value = 4578.8335
$ws.Range("I122").Value = 4896.6
$ws.Range("J122").Value = 2990
$ws.Range("K122").Value = 14689.8
$ws.Range("L122").Value = 8970
$ws.Range("M122").Value = -12239.8
$ws.Range("N122").Value = -13870
$ws.Range("H132").Value = 4014.182
$ws.Range("I132").Value = 3878.3794
$ws.Range("K132").Value = 11635.1382
$ws.Range("M132").Value = -9105.138199999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 349.6
$ws.Range("I9").Value = 399.5
$ws.Range("J9").Value = 150
$ws.Range("K9").Value = 399.5
$ws.Range("L9").Value = 150
$ws.Range("M9").Value = -175.5
$ws.Range("N9").Value = -598
$ws.Range("H17").Value = 11029
$ws.Range("I17").Value = 13038.667
$ws.Range("J17").Value = 5000
$ws.Range("K17").Value = 13038.667
$ws.Range("L17").Value = 5000
$ws.Range("M17").Value = -12868.667
$ws.Range("N17").Value = -5340
$ws.Range("H61").Value = 1710.6875
$ws.Range("I61").Value = 1891.2307
$ws.Range("K61").Value = 1891.2307
$ws.Range("M61").Value = -1689.2307
$ws.Range("H63").Value = 64999.668
$ws.Range("I63").Value = 29999.5
$ws.Range("J63").Value = 135000
$ws.Range("K63").Value = 29999.5
$ws.Range("L63").Value = 135000
$ws.Range("M63").Value = -29250.5
$ws.Range("N63").Value = -136498
$ws.Range("H66").Value = 64999.668
$ws.Range("I66").Value = 29999.5
$ws.Range("J66").Value = 135000
$ws.Range("K66").Value = 89998.5
$ws.Range("L66").Value = 405000
$ws.Range("M66").Value = -86254.5
$ws.Range("N66").Value = -412488
$ws.Range("H100").Value = 3639.1875
$ws.Range("I100").Value = 2714.4443
$ws.Range("J100").Value = 4828.143
$ws.Range("K100").Value = 2714.4443
$ws.Range("L100").Value = 4828.143
$ws.Range("M100").Value = -2173.4443
$ws.Range("N100").Value = -5910.143
$ws.Range("H113").Value = 1710.6875
$ws.Range("I113").Value = 1891.2307
$ws.Range("K113").Value = 1891.2307
$ws.Range("M113").Value = 278.7692999999999
$ws.Range("H122").Value = 8423.277
$ws.Range("J122").Value = 8093
$ws.Range("L122").Value = 24279
$ws.Range("N122").Value = -29179
$ws.Range("H132").Value = 43967.1
$ws.Range("I132").Value = 51346.25
$ws.Range("J132").Value = 8547.200000000001
$ws.Range("K132").Value = 154038.75
$ws.Range("L132").Value = 25641.6
$ws.Range("M132").Value = -151508.75
$ws.Range("N132").Value = -30701.6
$ws.Range("H136").Value = 3210.9
$ws.Range("I136").Value = 2849.5
$ws.Range("J136").Value = 3301.25
$ws.Range("K136").Value = 8548.5
$ws.Range("L136").Value = 9903.75
$ws.Range("M136").Value = -5998.5
$ws.Range("N136").Value = -15003.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H7").Value = 20004
$ws.Range("I7").Value = 20004
$ws.Range("K7").Value = 20004
$ws.Range("M7").Value = -19891
$ws.Range("H99").Value = 134000
$ws.Range("J99").Value = 134000
$ws.Range("L99").Value = 134000
$ws.Range("N99").Value = -139990
$ws.Range("H100").Value = 918
$ws.Range("I100").Value = 1050.5555
$ws.Range("K100").Value = 2101.111
$ws.Range("M100").Value = -1560.111
$ws.Range("H132").Value = 3103.0815
$ws.Range("I132").Value = 2460.1052
$ws.Range("K132").Value = 7380.3156
$ws.Range("M132").Value = -4850.3156
$ws.Range("H136").Value = 42469.36
$ws.Range("I136").Value = 2085.2144
$ws.Range("K136").Value = 6255.6432
$ws.Range("M136").Value = -3705.6432
